$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant (copy formatting only, reuse existing style records
# instead of minting new ones for every touched cell).
$xlPasteFormats = -4122

# Cells C5:F28 (rows 5, 12, 28) currently carry the "graded" fill styles
# (green/light-green). They should become plain bordered cells (same look
# as the already-blank G column in those rows) while getting a score of 5.
$plainFormatDonor = $ws.Range("G5")
$plainFormatDonor.Copy()

$rangesToClear = @("C5:F5", "C12:F12", "C28:F28", "C32:F32")
foreach ($addr in $rangesToClear) {
    $target = $ws.Range($addr)
    $target.PasteSpecial($xlPasteFormats)
    $target.Value2 = 5
}

$excel.CutCopyMode = $false

# Newly-scored single cells that already had the correct (unfilled) style —
# only the value needs to be set.
$ws.Range("G29").Value2 = 5
$ws.Range("G32").Value2 = 5

# Restore the active selection to E5 (bottom-right frozen pane).
$ws.Range("E5").Select()
